# Simplify image processing memory leaks
# Append one new data row (row 32) to each of the four log sheets,
# mirroring the existing row layout (time, lengths, checksum, decoded values).

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"
$newRow = 32
$timeValue = 45818.43747685185

function Set-NewRow($ws, $b, $c, $d, $e, $f, $g, $h, $i) {
    $ws.Cells.Item($newRow, 1).Value = $timeValue
    $ws.Cells.Item($newRow, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($newRow, 2).Value = $b
    $ws.Cells.Item($newRow, 3).Value = $c
    $ws.Cells.Item($newRow, 4).Value = $d
    $ws.Cells.Item($newRow, 5).Value = $e

    $ws.Cells.Item($newRow, 6).Value = $f
    $ws.Cells.Item($newRow, 7).Value = [double]$g
    $ws.Cells.Item($newRow, 8).Value = $h
    $ws.Cells.Item($newRow, 9).Value = $i
}

# DE_LFT_#1
$ws1 = $wb.Worksheets.Item("DE_LFT_#1")
Set-NewRow `
    $ws1 `
    "0x01,0x7c" `
    "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," `
    "0x01,0x74" `
    "0x14" `
    380 `
    "7.598631275147109e+23" `
    372 `
    14

# DE_LFT_#2
$ws2 = $wb.Worksheets.Item("DE_LFT_#2")
Set-NewRow `
    $ws2 `
    "0x01,0x7c" `
    "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," `
    "0x01,0x74" `
    "0xe" `
    380 `
    "5.68432987514711e+23" `
    372 `
    14

# DE_PLT_#1
$ws3 = $wb.Worksheets.Item("DE_PLT_#1")
Set-NewRow `
    $ws3 `
    "0x00,0x82" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," `
    "0x00,0x82" `
    "0x7" `
    130 `
    "5.68631262647114e+23" `
    129 `
    7

# DE_PLT_#2
$ws4 = $wb.Worksheets.Item("DE_PLT_#2")
Set-NewRow `
    $ws4 `
    "0x00,0x82" `
    "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," `
    "0x00,0x81" `
    "0x3" `
    130 `
    "9.85046333984776e+23" `
    129 `
    3
